# Scheduled runner update: refresh cached market-board pricing figures
# (currentAveragePrice / LevePrice / LeveProfit columns) across the Leve
# profit sheets, per the upstream data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 17166.666
$ws.Range("I18").Value = 17166.666
$ws.Range("K18").Value = 17166.666
$ws.Range("M18").Value = -16882.666
$ws.Range("H121").Value = 3256.2
$ws.Range("J121").Value = 3746
$ws.Range("L121").Value = 11238
$ws.Range("N121").Value = -14732
$ws.Range("H137").Value = 1995.1923
$ws.Range("I137").Value = 1873.8
$ws.Range("J137").Value = 2399.8333
$ws.Range("K137").Value = 5621.4
$ws.Range("L137").Value = 7199.499899999999
$ws.Range("M137").Value = -3071.4
$ws.Range("N137").Value = -12299.4999
$ws.Range("H138").Value = 2815.7896
$ws.Range("I138").Value = 2635
$ws.Range("J138").Value = 2880.3572
$ws.Range("K138").Value = 7905
$ws.Range("L138").Value = 8641.071599999999
$ws.Range("M138").Value = -2765
$ws.Range("N138").Value = -18921.0716

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 47998.332
$ws.Range("J23").Value = 46997.5
$ws.Range("L23").Value = 46997.5
$ws.Range("N23").Value = -47515.5
$ws.Range("H61").Value = 10717.228
$ws.Range("I61").Value = 10222.294
$ws.Range("K61").Value = 10222.294
$ws.Range("M61").Value = -10010.294
$ws.Range("H63").Value = 2200
$ws.Range("I63").Value = 2200
$ws.Range("K63").Value = 2200
$ws.Range("M63").Value = -1514
$ws.Range("H66").Value = 2200
$ws.Range("I66").Value = 2200
$ws.Range("K66").Value = 11000
$ws.Range("M66").Value = -7568
$ws.Range("H74").Value = 3923.1785
$ws.Range("J74").Value = 5599.8335
$ws.Range("L74").Value = 5599.8335
$ws.Range("N74").Value = -7347.8335
$ws.Range("H77").Value = 3923.1785
$ws.Range("J77").Value = 5599.8335
$ws.Range("L77").Value = 27999.1675
$ws.Range("N77").Value = -36735.1675
$ws.Range("H136").Value = 10717.228
$ws.Range("I136").Value = 10222.294
$ws.Range("K136").Value = 30666.882
$ws.Range("M136").Value = -28116.882

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 6338.4707
$ws.Range("I20").Value = 5429
$ws.Range("J20").Value = 7361.625
$ws.Range("K20").Value = 5429
$ws.Range("L20").Value = 7361.625
$ws.Range("M20").Value = -5182
$ws.Range("N20").Value = -7855.625
$ws.Range("H35").Value = 25800
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H94").Value = 967.51514
$ws.Range("I94").Value = 884.3333
$ws.Range("J94").Value = 1799.3334
$ws.Range("K94").Value = 884.3333
$ws.Range("L94").Value = 1799.3334
$ws.Range("M94").Value = -433.3333
$ws.Range("N94").Value = -2701.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 2795
$ws.Range("I15").Value = 1998.8889
$ws.Range("K15").Value = 1998.8889
$ws.Range("M15").Value = -1828.8889
$ws.Range("H22").Value = 698.375
$ws.Range("J22").Value = 1282.6666
$ws.Range("L22").Value = 1282.6666
$ws.Range("N22").Value = -1982.6666
$ws.Range("H60").Value = 14533
$ws.Range("I60").Value = 7549.5
$ws.Range("K60").Value = 7549.5
$ws.Range("M60").Value = -7038.5
$ws.Range("H68").Value = 39666.668
$ws.Range("I68").Value = 37000
$ws.Range("K68").Value = 37000
$ws.Range("M68").Value = -36251
$ws.Range("H71").Value = 39666.668
$ws.Range("I71").Value = 37000
$ws.Range("K71").Value = 111000
$ws.Range("M71").Value = -107256
$ws.Range("H74").Value = 42427
$ws.Range("J74").Value = 42427
$ws.Range("L74").Value = 42427
$ws.Range("N74").Value = -44175
$ws.Range("H77").Value = 42427
$ws.Range("J77").Value = 42427
$ws.Range("L77").Value = 127281
$ws.Range("N77").Value = -136017
$ws.Range("H132").Value = 3187
$ws.Range("I132").Value = 3108.1738
$ws.Range("K132").Value = 9324.5214
$ws.Range("M132").Value = -6794.5214

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 723.1111
$ws.Range("I8").Value = 723.1111
$ws.Range("K8").Value = 2169.3333
$ws.Range("M8").Value = -2030.3333
$ws.Range("H23").Value = 193.83333
$ws.Range("I23").Value = 192.33333
$ws.Range("J23").Value = 195.33333
$ws.Range("K23").Value = 576.99999
$ws.Range("L23").Value = 585.99999
$ws.Range("M23").Value = -341.99999
$ws.Range("N23").Value = -1055.99999
$ws.Range("H44").Value = 928665.0600000001
$ws.Range("J44").Value = 97.11539
$ws.Range("L44").Value = 291.34617
$ws.Range("N44").Value = -1087.34617
$ws.Range("H131").Value = 22729642
$ws.Range("J131").Value = 2598.7368
$ws.Range("L131").Value = 7796.2104
$ws.Range("N131").Value = -17876.2104

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 51604
$ws.Range("J46").Value = 51604
$ws.Range("L46").Value = 51604
$ws.Range("N46").Value = -51916
$ws.Range("H80").Value = 6930.1665
$ws.Range("J80").Value = 6050.8335
$ws.Range("L80").Value = 6050.8335
$ws.Range("N80").Value = -8046.8335
$ws.Range("H83").Value = 6930.1665
$ws.Range("J83").Value = 6050.8335
$ws.Range("L83").Value = 30254.1675
$ws.Range("N83").Value = -40238.1675
$ws.Range("H126").Value = 5240.6665
$ws.Range("I126").Value = 3748.1667
$ws.Range("J126").Value = 6733.1665
$ws.Range("K126").Value = 11244.5001
$ws.Range("L126").Value = 20199.4995
$ws.Range("M126").Value = -8774.500100000001
$ws.Range("N126").Value = -25139.4995
$ws.Range("H140").Value = 100000
$ws.Range("J140").Value = 100000
$ws.Range("L140").Value = 100000
$ws.Range("N140").Value = -110360

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 15000
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H22").Value = 1648.5714
$ws.Range("J22").Value = 1346.6666
$ws.Range("L22").Value = 1346.6666
$ws.Range("N22").Value = -1936.6666
$ws.Range("H27").Value = 1648.5714
$ws.Range("J27").Value = 1346.6666
$ws.Range("L27").Value = 1346.6666
$ws.Range("N27").Value = -1560.6666
$ws.Range("H46").Value = 14250
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 14250
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 14250
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -14626
$ws.Range("H98").Value = 49990
$ws.Range("J98").Value = 49990
$ws.Range("L98").Value = 49990
$ws.Range("N98").Value = -55980
$ws.Range("H100").Value = 3125.5908
$ws.Range("J100").Value = 3975.3333
$ws.Range("L100").Value = 3975.3333
$ws.Range("N100").Value = -5057.3333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H100").Value = 683
$ws.Range("I100").Value = 681.6667
$ws.Range("J100").Value = 687
$ws.Range("K100").Value = 1363.3334
$ws.Range("L100").Value = 1374
$ws.Range("M100").Value = -822.3334
$ws.Range("N100").Value = -2456
$ws.Range("H132").Value = 7500
$ws.Range("I132").Value = 7500
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 22500
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -19970
$ws.Range("N132").ClearContents()
